# Pedido 691ee840d9cd7e0aaa545c84
# - Limpia las celdas de texto vacias sobrantes (F19, G19, L19) de la fila 19.
# - Anade una nueva fila 20 con los datos del nuevo pedido (Mayte Lopez).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Algunas columnas (p.ej. "Unidades ...") contienen numeros guardados
    # como texto en el resto de la hoja. Forzamos formato de texto antes de
    # escribir para que Excel no los reinterprete como valores numericos,
    # y luego devolvemos el estilo de celda a "Normal" para no dejar
    # formato de texto aplicado de forma permanente.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# La fila 19 tenia celdas de texto vacias en F/G/L que ya no deben existir.
$ws.Range("F19").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("L19").ClearContents()

# Nueva fila 20.
$ws.Range("A20").Value = 2743
$ws.Range("B20").Value = "Mayte López"
$ws.Range("C20").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D20").Value = "MODULO FV JA SOLAR 535WP BLACK FRAME BIFACIAL 120 CELDAS"
Set-TextValue $ws.Range("E20") "10"
$ws.Range("H20").Value = "GOODWE GW5000-ES-20 híbrido monofásico"
Set-TextValue $ws.Range("I20") "1"
$ws.Range("J20").Value = "GOODWE Batería Lynx Home U G3 5,12 kWh"
Set-TextValue $ws.Range("K20") "1"
$ws.Range("M20").Value = "Sí"
$ws.Range("N20").Value = "2025-09-25T07:50:43.054Z"
